$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price (column D) values are stored as text (they use "." as a
# thousands separator for large coins, e.g. "70.102.20"), so every write to
# column D forces text via NumberFormat "@" and then restores the default
# 'Normal' cell style so no stray formatting is left behind.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.102.20'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.11%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.438.51'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.44%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '584.47'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.44%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '178.11'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.74%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.601'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.02%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.430.84'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.39%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.06%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.206'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.23%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.588'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.27%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '48.83'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.77%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000287'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.66%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '691.38'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.97%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.987.83'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.38%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.69'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.93%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '70.032.52'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.99%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.428.92'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.26%  '

$ws.Range('E19').Value = '  +0.95%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.73'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.50%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.47'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.17%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.903'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.09%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.57'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.05%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '17.08'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.09%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '101.15'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.10%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.94'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.67%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.68'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.55%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.63'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.09%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.61'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.38%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.79'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.07%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.21'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.37%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.84'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.97%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '575.30'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.64%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.05'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.92%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '58.82'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.40%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.104'
$ws.Range('D36').Style = 'Normal'

$ws.Range('E37').Value = '  +0.05%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.593.12'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.53%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.140'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.44%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '35.39'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.24%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0₃0744'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.60%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.33'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.53%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.69'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.54%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.35'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.43%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0423'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.41%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.335'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.88%  '

$ws.Range('B47').Value = 'ThetaToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.69'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.76%  '

$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.45'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.32%  '

$ws.Range('E49').Value = '  +0.01%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.997'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.37%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '133.56'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.97%  '
